$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.445.78"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "3.692.54"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'686.40"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").Value = "'160.31"
$ws.Range("E6").Value = "  -5.61%  "
$ws.Range("D7").Value = "3.691.13"
$ws.Range("E7").Value = "  -2.89%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -5.90%  "
$ws.Range("D10").Value = "'0.147"
$ws.Range("E10").Value = "  -8.45%  "
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("D12").Value = "'0.436"
$ws.Range("E12").Value = "  -9.36%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -6.71%  "
$ws.Range("D14").Value = "4.313.99"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "'32.53"
$ws.Range("E15").Value = "  -10.21%  "
$ws.Range("D16").Value = "3.679.30"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "69.450.04"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -9.51%  "
$ws.Range("E20").Value = "  -10.62%  "
$ws.Range("D21").Value = "'470.96"
$ws.Range("E21").Value = "  -8.21%  "
$ws.Range("D22").Value = "'9.99"
$ws.Range("E22").Value = "  -4.57%  "
$ws.Range("D23").Value = "'0.649"
$ws.Range("E23").Value = "  -9.43%  "
$ws.Range("D24").Value = "'79.69"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("D25").Value = "3.837.85"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("D27").Value = "'0.0000126"
$ws.Range("E27").Value = "  -10.98%  "
$ws.Range("E28").Value = "  -13.06%  "
$ws.Range("E29").Value = "  -10.26%  "
$ws.Range("E30").Value = "  -9.45%  "
$ws.Range("E31").Value = "  -12.30%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.02"
$ws.Range("E32").Value = "  -10.88%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.66"
$ws.Range("E33").Value = "  -8.65%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'26.79"
$ws.Range("E35").Value = "  -8.06%  "
$ws.Range("E36").Value = "  -6.85%  "
$ws.Range("D37").Value = "'8.21"
$ws.Range("E37").Value = "  -11.91%  "
$ws.Range("D38").Value = "'6.18"
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("E39").Value = "  -2.99%  "
$ws.Range("D41").Value = "'0.0906"
$ws.Range("E41").Value = "  -10.10%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'167.39"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.944"
$ws.Range("E44").Value = "  -6.66%  "
$ws.Range("D45").Value = "'47.92"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "  -13.81%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'28.68"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.31"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "'0.000277"
$ws.Range("E50").Value = "  -9.17%  "
$ws.Range("D51").Value = "'375.09"
$ws.Range("E51").Value = "  -12.80%  "
